# Generate Report for Handoff
# Updates the status and timestamps on the Overview, zh-cn, and de-de sheets
# to reflect the files being ready for / at handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (B2, C2) and Latest Handoff Date (D2)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-31-18 20:31:41"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-18 20:31:38"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-18 20:31:41"
